$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1, G1, H1 - copy style from the existing header (E1) and set text
$headers = @("KNN_Outliers_MAD", "SVM_Outliers_MAD", "RF_Outliers_MAD")
$cols = @("F", "G", "H")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $cell = $ws.Range($col + "1")
    $cell.Value = $headers[$i]
    $ws.Range("E1").Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
    $cell.Value = $headers[$i]
}

# Fill F2:H16 with boolean FALSE values
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("F" + $r).Value = $false
    $ws.Range("G" + $r).Value = $false
    $ws.Range("H" + $r).Value = $false
}
